# TC2 - invalid login test data
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# containing a UserName/Password header row and an invalid bhanu/bhanu123
# credential row, then restores the selection on "ValidLogin" and makes the
# new sheet the active tab (as it was left selected last in Excel).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new sheet right after "ValidLogin" ---------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "InvalidLogin"

# --- Populate the invalid-login test data -------------------------------
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "bhanu"
$ws2.Range("B2").Value = "bhanu123"

# --- Restore the ValidLogin sheet's selection to A1:B2 ------------------
[void]$ws1.Select()
[void]$ws1.Range("A1:B2").Select()

# --- Select InvalidLogin (B3) and bump its zoom; leave it as active tab -
[void]$ws2.Select()
[void]$ws2.Range("B3").Select()
$ws2.Application.ActiveWindow.Zoom = 160
